# Fill in the "22-10-2021"/"23-10-2021"/"25-10-2021" start/target/finish dates
# for the "Analis Dataset" (row 3), "Modelling menggunakan Algoritma CNN Part 1"
# (row 5) and "Evaluasi Model CNN 1" (row 6) tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C5 has no explicit style yet (unlike the other date cells in column C/D/E,
# which already carry the bordered "text" number format). Copy that look from
# a neighbouring cell that already has it before writing the C5 value, so the
# new cell matches the rest of the column.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 3 - "Analis Dataset"
$ws.Range("C3").Value = "22-10-2021"
$ws.Range("D3").Value = "23-10-2021"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "23-10-2021"

# Row 5 - "Modelling menggunakan Algoritma CNN Part 1"
$ws.Range("C5").Value = "25-10-2021"
$ws.Range("D5").Value = "25-10-2021"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "25-10-2021"

# Row 6 - "Evaluasi Model CNN 1"
$ws.Range("C6").Value = "25-10-2021"
$ws.Range("D6").Value = "25-10-2021"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "25-10-2021"

# Match the author's active selection when they saved the file.
$ws.Range("E5").Select() | Out-Null
